# Auto-generated edit script: updates the Cryptos sheet per the target diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price column (D): preserve original text type/style by reading the
# existing Style before writing, prefixing the literal with an apostrophe
# so numeric-looking strings (e.g. "570.47") are stored as text (matching
# the source workbook, where every Price cell is a string), then restoring
# the original Style so no stray number-format/quote-prefix style sticks.
$origStyle = $ws.Range("D2").Style
$ws.Range("D2").Value = "'60.758.48"
$ws.Range("D2").Style = $origStyle
$origStyle = $ws.Range("D3").Style
$ws.Range("D3").Value = "'2.399.47"
$ws.Range("D3").Style = $origStyle
$origStyle = $ws.Range("D5").Style
$ws.Range("D5").Value = "'570.47"
$ws.Range("D5").Style = $origStyle
$origStyle = $ws.Range("D6").Style
$ws.Range("D6").Value = "'139.70"
$ws.Range("D6").Style = $origStyle
$origStyle = $ws.Range("D8").Style
$ws.Range("D8").Value = "'0.526"
$ws.Range("D8").Style = $origStyle
$origStyle = $ws.Range("D9").Style
$ws.Range("D9").Value = "'2.379.27"
$ws.Range("D9").Style = $origStyle
$origStyle = $ws.Range("D12").Style
$ws.Range("D12").Value = "'5.07"
$ws.Range("D12").Style = $origStyle
$origStyle = $ws.Range("D13").Style
$ws.Range("D13").Value = "'0.338"
$ws.Range("D13").Style = $origStyle
$origStyle = $ws.Range("D14").Style
$ws.Range("D14").Value = "'25.92"
$ws.Range("D14").Style = $origStyle
$origStyle = $ws.Range("D15").Style
$ws.Range("D15").Value = "'0.0000170"
$ws.Range("D15").Style = $origStyle
$origStyle = $ws.Range("D17").Style
$ws.Range("D17").Value = "'60.775.64"
$ws.Range("D17").Style = $origStyle
$origStyle = $ws.Range("D18").Style
$ws.Range("D18").Value = "'2.382.38"
$ws.Range("D18").Style = $origStyle
$origStyle = $ws.Range("D19").Style
$ws.Range("D19").Value = "'10.52"
$ws.Range("D19").Style = $origStyle
$origStyle = $ws.Range("D20").Style
$ws.Range("D20").Value = "'7.19"
$ws.Range("D20").Style = $origStyle
$origStyle = $ws.Range("D21").Style
$ws.Range("D21").Value = "'321.80"
$ws.Range("D21").Style = $origStyle
$origStyle = $ws.Range("D22").Style
$ws.Range("D22").Value = "'4.02"
$ws.Range("D22").Style = $origStyle
$origStyle = $ws.Range("D26").Style
$ws.Range("D26").Value = "'64.27"
$ws.Range("D26").Style = $origStyle
$origStyle = $ws.Range("D27").Style
$ws.Range("D27").Value = "'8.61"
$ws.Range("D27").Style = $origStyle
$origStyle = $ws.Range("D28").Style
$ws.Range("D28").Value = "'570.95"
$ws.Range("D28").Style = $origStyle
$origStyle = $ws.Range("D29").Style
$ws.Range("D29").Value = "'2.514.99"
$ws.Range("D29").Style = $origStyle
$origStyle = $ws.Range("D30").Style
$ws.Range("D30").Value = "'0.0₃0910"
$ws.Range("D30").Style = $origStyle
$origStyle = $ws.Range("D31").Style
$ws.Range("D31").Value = "'7.81"
$ws.Range("D31").Style = $origStyle
$origStyle = $ws.Range("D32").Style
$ws.Range("D32").Value = "'1.35"
$ws.Range("D32").Style = $origStyle
$origStyle = $ws.Range("D36").Style
$ws.Range("D36").Value = "'4.61"
$ws.Range("D36").Style = $origStyle
$origStyle = $ws.Range("D37").Style
$ws.Range("D37").Value = "'0.367"
$ws.Range("D37").Style = $origStyle
$origStyle = $ws.Range("D38").Style
$ws.Range("D38").Value = "'1.38"
$ws.Range("D38").Style = $origStyle
$origStyle = $ws.Range("D39").Style
$ws.Range("D39").Value = "'18.17"
$ws.Range("D39").Style = $origStyle
$origStyle = $ws.Range("D40").Style
$ws.Range("D40").Value = "'146.70"
$ws.Range("D40").Style = $origStyle
$origStyle = $ws.Range("D41").Style
$ws.Range("D41").Value = "'5.07"
$ws.Range("D41").Style = $origStyle
$origStyle = $ws.Range("D43").Style
$ws.Range("D43").Value = "'41.51"
$ws.Range("D43").Style = $origStyle
$origStyle = $ws.Range("D44").Style
$ws.Range("D44").Value = "'1.66"
$ws.Range("D44").Style = $origStyle
$origStyle = $ws.Range("D45").Style
$ws.Range("D45").Value = "'2.33"
$ws.Range("D45").Style = $origStyle
$origStyle = $ws.Range("D47").Style
$ws.Range("D47").Value = "'140.02"
$ws.Range("D47").Style = $origStyle
$origStyle = $ws.Range("D48").Style
$ws.Range("D48").Value = "'3.50"
$ws.Range("D48").Style = $origStyle
$origStyle = $ws.Range("D51").Style
$ws.Range("D51").Value = "'19.26"
$ws.Range("D51").Style = $origStyle

# Coin / Link / Volume(1h) columns: plain text assignment (never
# auto-coerced to numbers by Excel, since they contain letters, URLs, or
# padding spaces + a trailing "%").
$ws.Range("E2").Value = "  -2.30%  "
$ws.Range("E3").Value = "  -2.15%  "
$ws.Range("E4").Value = "  -0.32%  "
$ws.Range("E5").Value = "  -1.67%  "
$ws.Range("E6").Value = "  -2.70%  "
$ws.Range("E7").Value = "  +0.23%  "
$ws.Range("E8").Value = "  -0.92%  "
$ws.Range("E9").Value = "  -2.87%  "
$ws.Range("E10").Value = "  -0.09%  "
$ws.Range("E11").Value = "  +0.30%  "
$ws.Range("E12").Value = "  -2.61%  "
$ws.Range("E13").Value = "  -1.80%  "
$ws.Range("E14").Value = "  -2.25%  "
$ws.Range("E15").Value = "  -1.85%  "
$ws.Range("E16").Value = "  -0.32%  "
$ws.Range("E17").Value = "  -2.12%  "
$ws.Range("E18").Value = "  -2.12%  "
$ws.Range("E19").Value = "  -3.21%  "
$ws.Range("E20").Value = "  +0.50%  "
$ws.Range("E21").Value = "  -2.34%  "
$ws.Range("E22").Value = "  -1.95%  "
$ws.Range("E23").Value = "  +1.67%  "
$ws.Range("E24").Value = "  +0.17%  "
$ws.Range("E25").Value = "  -6.69%  "
$ws.Range("E26").Value = "  -2.51%  "
$ws.Range("E27").Value = "  -8.07%  "
$ws.Range("E28").Value = "  -7.72%  "
$ws.Range("E29").Value = "  -1.04%  "
$ws.Range("E30").Value = "  -4.98%  "
$ws.Range("E31").Value = "  -2.37%  "
$ws.Range("E32").Value = "  -6.07%  "
$ws.Range("E33").Value = "  -2.22%  "
$ws.Range("E34").Value = "  -7.13%  "
$ws.Range("E35").Value = "  +0.30%  "
$ws.Range("E36").Value = "  -5.85%  "
$ws.Range("E37").Value = "  -2.53%  "
$ws.Range("E38").Value = "  -3.68%  "
$ws.Range("B39").Value = "EthereumClassic"
$ws.Range("C39").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("E39").Value = "  -1.04%  "
$ws.Range("B40").Value = "Monero"
$ws.Range("C40").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("E40").Value = "  -2.39%  "
$ws.Range("E41").Value = "  -4.73%  "
$ws.Range("E42").Value = "  +0.13%  "
$ws.Range("E43").Value = "  -2.36%  "
$ws.Range("E44").Value = "  -4.90%  "
$ws.Range("E45").Value = "  -5.34%  "
$ws.Range("E46").Value = "  +17.01%  "
$ws.Range("E47").Value = "  -2.25%  "
$ws.Range("E48").Value = "  -3.90%  "
$ws.Range("E49").Value = "  -3.47%  "
$ws.Range("E50").Value = "  -4.22%  "
$ws.Range("E51").Value = "  -1.55%  "
